# 자동 업데이트: 2025-04-12 12:45:01
# Remove the "-100" offset from the percentage-change formulas in columns
# N (vs. day-1 baseline D2) and Q (vs. day-1 baseline O2), so that the
# columns report "index = 100" on day 1 instead of "change = 0%".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column N: row 2 has its own (non-shared) formula; rows 3:41 are one
# shared-formula group anchored at N3.
$ws.Range("N2").Formula = '=M2/$D$2*100'
$ws.Range("N3:N41").Formula = '=M3/$D$2*100'

# Column Q: Q2 is a literal value (index base = 100, was 0); row 3 has its
# own formula; rows 4:41 are a shared-formula group anchored at Q4.
$ws.Range("Q2").Value = 100
$ws.Range("Q3").Formula = '=P3/$O$2*100'
$ws.Range("Q4:Q41").Formula = '=P4/$O$2*100'

# Update the saved selection to match the author's final cursor position.
$ws.Range("Q3:Q41").Select()
